$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.016.56'
$ws.Cells.Item(2, 5).Value = '  +0.34%  '
$ws.Cells.Item(3, 4).Value = '3.108.98'
$ws.Cells.Item(3, 5).Value = '  +0.38%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Value = '578.41'
$ws.Cells.Item(5, 5).Value = '  -0.14%  '
$ws.Cells.Item(6, 4).Value = '173.27'
$ws.Cells.Item(6, 5).Value = '  +0.42%  '
$ws.Cells.Item(7, 5).Value = '  +0.12%  '
$ws.Cells.Item(8, 4).Value = '0.520'
$ws.Cells.Item(8, 5).Value = '  -0.40%  '
$ws.Cells.Item(9, 4).Value = '6.48'
$ws.Cells.Item(9, 5).Value = '  +0.65%  '
$ws.Cells.Item(11, 5).Value = '  -0.78%  '
$ws.Cells.Item(12, 5).Value = '  -0.42%  '
$ws.Cells.Item(13, 4).Value = '36.78'
$ws.Cells.Item(13, 5).Value = '  -1.60%  '
$ws.Cells.Item(14, 5).Value = '  -1.61%  '
$ws.Cells.Item(15, 4).Value = '3.625.65'
$ws.Cells.Item(15, 5).Value = '  +0.48%  '
$ws.Cells.Item(16, 4).Value = '67.000.63'
$ws.Cells.Item(16, 5).Value = '  +0.40%  '
$ws.Cells.Item(17, 5).Value = '  -1.36%  '
$ws.Cells.Item(18, 4).Value = '3.111.90'
$ws.Cells.Item(18, 5).Value = '  +0.51%  '
$ws.Cells.Item(19, 4).Value = '16.55'
$ws.Cells.Item(19, 5).Value = '  +1.63%  '
$ws.Cells.Item(20, 4).Value = '491.41'
$ws.Cells.Item(20, 5).Value = '  +2.22%  '
$ws.Cells.Item(21, 5).Value = '  -1.49%  '
$ws.Cells.Item(22, 5).Value = '  +4.40%  '
$ws.Cells.Item(23, 4).Value = '83.96'
$ws.Cells.Item(23, 5).Value = '  -0.08%  '
$ws.Cells.Item(24, 4).Value = '13.09'
$ws.Cells.Item(24, 5).Value = '  -1.04%  '
$ws.Cells.Item(25, 5).Value = '  -2.95%  '
$ws.Cells.Item(26, 5).Value = '  +4.18%  '
$ws.Cells.Item(27, 5).Value = '  -0.03%  '
$ws.Cells.Item(28, 4).Value = '7.90'
$ws.Cells.Item(28, 5).Value = '  -1.18%  '
$ws.Cells.Item(29, 5).Value = '  -1.60%  '
$ws.Cells.Item(30, 5).Value = '  -0.70%  '
$ws.Cells.Item(31, 4).Value = '28.36'
$ws.Cells.Item(31, 5).Value = '  -1.53%  '
$ws.Cells.Item(32, 5).Value = '  -1.29%  '
$ws.Cells.Item(33, 4).Value = '0.0₃0944'
$ws.Cells.Item(33, 5).Value = '  -6.07%  '
$ws.Cells.Item(34, 4).Value = '0.999'
$ws.Cells.Item(34, 5).Value = '  -0.03%  '
$ws.Cells.Item(35, 5).Value = '  -0.52%  '
$ws.Cells.Item(36, 4).Value = '0.971'
$ws.Cells.Item(36, 5).Value = '  -2.15%  '
$ws.Cells.Item(37, 4).Value = '47.34'
$ws.Cells.Item(37, 5).Value = '  -1.26%  '
$ws.Cells.Item(38, 5).Value = '  -3.23%  '
$ws.Cells.Item(39, 5).Value = '  -2.97%  '
$ws.Cells.Item(40, 5).Value = '  +1.13%  '
$ws.Cells.Item(41, 5).Value = '  -2.37%  '
$ws.Cells.Item(42, 4).Value = '385.42'
$ws.Cells.Item(42, 5).Value = '  +0.42%  '
$ws.Cells.Item(43, 4).Value = '2.811.35'
$ws.Cells.Item(43, 5).Value = '  -1.06%  '
$ws.Cells.Item(44, 4).Value = '2.60'
$ws.Cells.Item(44, 5).Value = '  -7.43%  '
$ws.Cells.Item(45, 5).Value = '  -2.50%  '
$ws.Cells.Item(46, 4).Value = '135.42'
$ws.Cells.Item(46, 5).Value = '  +0.06%  '
$ws.Cells.Item(48, 4).Value = '24.82'
$ws.Cells.Item(48, 5).Value = '  -0.86%  '
$ws.Cells.Item(49, 5).Value = '  -1.61%  '
$ws.Cells.Item(50, 5).Value = '  -1.06%  '
$ws.Cells.Item(51, 4).Value = '6.71'
$ws.Cells.Item(51, 5).Value = '  -2.04%  '
